$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped entirely: "RM 232" (row 26) and "SC 92" (row 28, which becomes row 27
# after the first deletion shifts everything up).
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# Apply the remaining per-cell value changes (after the row shift) to match the updated imputation results.
$ws.Cells.Item(2, 4).Value = ""
$ws.Cells.Item(5, 4).Value = -14.4
$ws.Cells.Item(6, 3).Value = 15.1
$ws.Cells.Item(6, 4).Value = -14.2
$ws.Cells.Item(8, 3).Value = ""
$ws.Cells.Item(10, 4).Value = ""
$ws.Cells.Item(12, 3).Value = 12.5
$ws.Cells.Item(13, 4).Value = ""
$ws.Cells.Item(14, 3).Value = ""
$ws.Cells.Item(17, 3).Value = 11.2
$ws.Cells.Item(18, 3).Value = 11.5
$ws.Cells.Item(19, 3).Value = ""
$ws.Cells.Item(20, 3).Value = ""
$ws.Cells.Item(23, 3).Value = 12.2
$ws.Cells.Item(24, 4).Value = -13.9
$ws.Cells.Item(27, 2).Value = -20.4
$ws.Cells.Item(27, 3).Value = ""
$ws.Cells.Item(28, 4).Value = ""
$ws.Cells.Item(29, 2).Value = ""
$ws.Cells.Item(30, 4).Value = -13.6
$ws.Cells.Item(32, 2).Value = ""
